$p = $ppt.ActivePresentation

# --- 1) Update the auto-date placeholder text on the slide master and
#        every slide layout (the "datetimeFigureOut" field caches its
#        last-rendered string in <a:t>; PowerPoint refreshes that cache on
#        save -- here we just rewrite the cached text directly). ---
$newDate = "8/29/2017"

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.Name -like "Date Placeholder*") {
        $shape.TextFrame.TextRange.Text = $newDate
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shape = $layout.Shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- 2) Bump the revision number shown on the binder-cover slide. ---
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTextFrame -and $shape.TextFrame.TextRange.Text -like "Version *") {
        $shape.TextFrame.TextRange.Text = "Version 1.4"
    }
}
